$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 737.1539
$ws.Range("I33").Value = 816.63635
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 816.63635
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -587.63635
$ws.Range("N33").Value = -758

$ws.Range("H53").Value = 222
$ws.Range("I53").Value = 215.2
$ws.Range("J53").Value = 233.33333
$ws.Range("K53").Value = 215.2
$ws.Range("L53").Value = 233.33333
$ws.Range("M53").Value = 421.8
$ws.Range("N53").Value = -1507.33333

$ws.Range("H64").Value = 2958.3333
$ws.Range("I64").Value = 2750
$ws.Range("K64").Value = 2750
$ws.Range("M64").Value = -2502

$ws.Range("H67").Value = 2958.3333
$ws.Range("I67").Value = 2750
$ws.Range("K67").Value = 2750
$ws.Range("M67").Value = -1892

$ws.Range("H121").Value = 1895
$ws.Range("J121").Value = 2300.625
$ws.Range("L121").Value = 6901.875
$ws.Range("N121").Value = -10395.875

$ws.Range("H132").Value = 1830.2642
$ws.Range("I132").Value = 1841.2549
$ws.Range("J132").Value = 1550
$ws.Range("K132").Value = 5523.7647
$ws.Range("L132").Value = 4650
$ws.Range("M132").Value = -2993.7647
$ws.Range("N132").Value = -9710

$ws.Range("H135").Value = 21740338
$ws.Range("I135").Value = 26316888
$ws.Range("J135").Value = 1725
$ws.Range("K135").Value = 236851992
$ws.Range("L135").Value = 15525
$ws.Range("M135").Value = -236849457
$ws.Range("N135").Value = -20595

$ws.Range("H137").Value = 2138194.8
$ws.Range("I137").Value = 4630733.5
$ws.Range("J137").Value = 1732.619
$ws.Range("K137").Value = 13892200.5
$ws.Range("L137").Value = 5197.857
$ws.Range("M137").Value = -13889650.5
$ws.Range("N137").Value = -10297.857

$ws.Range("H140").Value = 74548.336
$ws.Range("J140").Value = 74548.336
$ws.Range("L140").Value = 74548.336
$ws.Range("N140").Value = -84908.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6419619
$ws.Range("I32").Value = 7050962
$ws.Range("J32").Value = 16000
$ws.Range("K32").Value = 7050962
$ws.Range("L32").Value = 16000
$ws.Range("M32").Value = -7050675
$ws.Range("N32").Value = -16574

$ws.Range("H61").Value = 12824147
$ws.Range("I61").Value = 19609988
$ws.Range("J61").Value = 6444.4443
$ws.Range("K61").Value = 19609988
$ws.Range("L61").Value = 6444.4443
$ws.Range("M61").Value = -19609776
$ws.Range("N61").Value = -6868.4443

$ws.Range("H74").Value = 16131908
$ws.Range("I74").Value = 2287.3635
$ws.Range("J74").Value = 55559868
$ws.Range("K74").Value = 2287.3635
$ws.Range("L74").Value = 55559868
$ws.Range("M74").Value = -1413.3635
$ws.Range("N74").Value = -55561616

$ws.Range("H77").Value = 16131908
$ws.Range("I77").Value = 2287.3635
$ws.Range("J77").Value = 55559868
$ws.Range("K77").Value = 11436.8175
$ws.Range("L77").Value = 277799340
$ws.Range("M77").Value = -7068.817499999999
$ws.Range("N77").Value = -277808076

$ws.Range("H122").Value = 60499.883
$ws.Range("I122").Value = 84705.836
$ws.Range("J122").Value = 2405.6
$ws.Range("K122").Value = 254117.508
$ws.Range("L122").Value = 7216.799999999999
$ws.Range("M122").Value = -251667.508
$ws.Range("N122").Value = -12116.8

$ws.Range("H132").Value = 2657054.2
$ws.Range("I132").Value = 3624.6
$ws.Range("J132").Value = 8553565
$ws.Range("K132").Value = 10873.8
$ws.Range("L132").Value = 25660695
$ws.Range("M132").Value = -8343.799999999999
$ws.Range("N132").Value = -25665755

$ws.Range("H136").Value = 12824147
$ws.Range("I136").Value = 19609988
$ws.Range("J136").Value = 6444.4443
$ws.Range("K136").Value = 58829964
$ws.Range("L136").Value = 19333.3329
$ws.Range("M136").Value = -58827414
$ws.Range("N136").Value = -24433.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6555.949
$ws.Range("I31").Value = 1853.6
$ws.Range("J31").Value = 7515.6123
$ws.Range("K31").Value = 1853.6
$ws.Range("L31").Value = 7515.6123
$ws.Range("M31").Value = -1558.6
$ws.Range("N31").Value = -8105.6123

$ws.Range("H34").Value = 6555.949
$ws.Range("I34").Value = 1853.6
$ws.Range("J34").Value = 7515.6123
$ws.Range("K34").Value = 1853.6
$ws.Range("L34").Value = 7515.6123
$ws.Range("M34").Value = -1651.6
$ws.Range("N34").Value = -7919.6123

$ws.Range("H58").Value = 1576.75
$ws.Range("I58").Value = 1313.1666
$ws.Range("J58").Value = 1840.3334
$ws.Range("K58").Value = 1313.1666
$ws.Range("L58").Value = 1840.3334
$ws.Range("M58").Value = -1110.1666
$ws.Range("N58").Value = -2246.3334

$ws.Range("H68").Value = 23199.455
$ws.Range("J68").Value = 23199.455
$ws.Range("L68").Value = 23199.455
$ws.Range("N68").Value = -24697.455

$ws.Range("H71").Value = 23199.455
$ws.Range("J71").Value = 23199.455
$ws.Range("L71").Value = 69598.36500000001
$ws.Range("N71").Value = -77086.36500000001

$ws.Range("H74").Value = 19635.818
$ws.Range("J74").Value = 19635.818
$ws.Range("L74").Value = 19635.818
$ws.Range("N74").Value = -21383.818

$ws.Range("H77").Value = 19635.818
$ws.Range("J77").Value = 19635.818
$ws.Range("L77").Value = 58907.454
$ws.Range("N77").Value = -67643.454

$ws.Range("H132").Value = 25643542
$ws.Range("I132").Value = 26318048
$ws.Range("J132").Value = 23812738
$ws.Range("K132").Value = 78954144
$ws.Range("L132").Value = 71438214
$ws.Range("M132").Value = -78951614
$ws.Range("N132").Value = -71443274

$ws.Range("H136").Value = 1576.75
$ws.Range("I136").Value = 1313.1666
$ws.Range("J136").Value = 1840.3334
$ws.Range("K136").Value = 3939.4998
$ws.Range("L136").Value = 5521.0002
$ws.Range("M136").Value = -1389.4998
$ws.Range("N136").Value = -10621.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 7380.2856
$ws.Range("I68").Value = 11962.75
$ws.Range("K68").Value = 35888.25
$ws.Range("M68").Value = -35077.25

$ws.Range("H71").Value = 7380.2856
$ws.Range("I71").Value = 11962.75
$ws.Range("K71").Value = 107664.75
$ws.Range("M71").Value = -103608.75

$ws.Range("H113").Value = 690.2857
$ws.Range("I113").Value = 862.5
$ws.Range("K113").Value = 2587.5
$ws.Range("M113").Value = -417.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 5000
$ws.Range("J6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5226

$ws.Range("H16").Value = 5000
$ws.Range("J16").Value = 5000
$ws.Range("L16").Value = 5000
$ws.Range("N16").Value = -5500

$ws.Range("H132").Value = 33338878
$ws.Range("I132").Value = 40005492
$ws.Range("J132").Value = 5802.2
$ws.Range("K132").Value = 120016476
$ws.Range("L132").Value = 17406.6
$ws.Range("M132").Value = -120013946
$ws.Range("N132").Value = -22466.6

$ws.Range("H141").Value = 60932
$ws.Range("J141").Value = 60932
$ws.Range("L141").Value = 60932
$ws.Range("N141").Value = -71292

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11567
$ws.Range("I22").Value = 578.3333
$ws.Range("J22").Value = 33544.332
$ws.Range("K22").Value = 578.3333
$ws.Range("L22").Value = 33544.332
$ws.Range("M22").Value = -283.3333
$ws.Range("N22").Value = -34134.332

$ws.Range("H27").Value = 11567
$ws.Range("I27").Value = 578.3333
$ws.Range("J27").Value = 33544.332
$ws.Range("K27").Value = 578.3333
$ws.Range("L27").Value = 33544.332
$ws.Range("M27").Value = -471.3333
$ws.Range("N27").Value = -33758.332

$ws.Range("H82").Value = 3353.875
$ws.Range("I82").Value = 4666.3335
$ws.Range("J82").Value = 2566.4
$ws.Range("K82").Value = 4666.3335
$ws.Range("L82").Value = 2566.4
$ws.Range("M82").Value = -4305.3335
$ws.Range("N82").Value = -3288.4

$ws.Range("H85").Value = 3353.875
$ws.Range("I85").Value = 4666.3335
$ws.Range("J85").Value = 2566.4
$ws.Range("K85").Value = 4666.3335
$ws.Range("L85").Value = 2566.4
$ws.Range("M85").Value = -3418.3335
$ws.Range("N85").Value = -5062.4

$ws.Range("H136").Value = 1789.9333
$ws.Range("I136").Value = 1742.7858
$ws.Range("J136").Value = 2450
$ws.Range("K136").Value = 5228.357400000001
$ws.Range("L136").Value = 7350
$ws.Range("M136").Value = -2678.357400000001
$ws.Range("N136").Value = -12450

$ws.Range("H139").Value = 1176824.4
$ws.Range("J139").Value = 44125.938
$ws.Range("L139").Value = 44125.938
$ws.Range("N139").Value = -54405.938

$ws.Range("H140").Value = 64746.11
$ws.Range("J140").Value = 64746.11
$ws.Range("L140").Value = 64746.11
$ws.Range("N140").Value = -75106.11
